$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CMSgrantHistTable")

$ws.Range("A3").Value = "This table shows the grant awards and award dollars CMS made for FY 2012-2016. It is provided as a text alternative to the interactive chart on the CMS page of this website."
$ws.Range("A7").Value = "Grant awards and award dollars CMS made for FY 2012-2016."
